{"js": "// Insert a new bulleted list item (with a hyperlink) right after the\n// \"https://www2.cs.arizona.edu/people/mccann/errors-python#Four\" bullet\n// and right before the trailing empty bookmark paragraph.\n\nconst url = \"https://guides.github.com/features/mastering-markdown/\";\n\nconst body = context.document.body;\nbody.paragraphs.load(\"items\");\nawait context.sync();\n\nconst paragraphs = body.paragraphs.items;\n// Last paragraph is the empty one holding the _GoBack bookmark; the one\n// right before it is the \"...errors-python#Four\" bullet we anchor on.\nconst anchorParagraph = paragraphs[paragraphs.length - 2];\n\n// Inserting \"After\" the anchor paragraph creates a sibling paragraph that\n// inherits the same list formatting (ListParagraph style + numPr).\nconst newParagraph = anchorParagraph.insertParagraph(\"\", \"After\");\nawait context.sync();\n\n// Insert the URL text, then turn just that text into a hyperlink (doing it\n// via the returned Range \u2014 rather than the whole paragraph, which would\n// also include the paragraph mark \u2014 avoids a stray empty trailing run).\nconst textRange = newParagraph.insertText(url, \"Replace\");\ntextRange.hyperlink = url;\nawait context.sync();\n", "ps1": "# Insert a new bulleted list item (with a hyperlink) right after the\n# \"https://www2.cs.arizona.edu/people/mccann/errors-python#Four\" bullet\n# and right before the trailing empty bookmark paragraph.\n\n$d = $word.ActiveDocument\n\n$url = \"https://guides.github.com/features/mastering-markdown/\"\n\n# Locate the last paragraph (empty paragraph holding the _GoBack bookmark)\n# and the paragraph right before it (the \"...errors-python#Four\" bullet).\n$paraCount = $d.Paragraphs.Count\n$anchorPara = $d.Paragraphs.Item($paraCount - 1)\n\n# Insert a new paragraph right after the anchor paragraph; it inherits the\n# same list formatting (ListParagraph style + numPr) automatically.\n$anchorRange = $anchorPara.Range\n$anchorRange.Collapse(0)   # wdCollapseEnd\n$anchorRange.InsertParagraphAfter()\n\n# Re-fetch the freshly created (now previous-to-last) paragraph and give it\n# the link text.\n$paraCount = $d.Paragraphs.Count\n$newPara = $d.Paragraphs.Item($paraCount - 1)\n$newParaRange = $newPara.Range\n$newParaRange.Text = $url\n\n# Turn just the URL text (not the paragraph mark) into a hyperlink.\n$textRange = $d.Range($newParaRange.Start, $newParaRange.Start + $url.Length)\n$d.Hyperlinks.Add($textRange, $url)\n"}
